$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SignUp")

# Update the sign-up email addresses (govind / namrata got a digit appended)
# and the interests column for the Gajendra Rathod row.
$ws.Range("B2").Value = "govind1@asite.com"
$ws.Range("B3").Value = "namrata2@asite.com"
$ws.Range("D4").Value = "Healthcare"

# Turn the email addresses into real mailto hyperlinks.
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:govind1@asite.com")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:namrata2@asite.com")
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:gajendra@asite.com")

# Leave the selection where the author left off.
$ws.Range("D18").Select()
